# signal-tracker-v3 / signals.xlsx
#
# Commit: "delete pip installs and change sheet name to primary"
#
# The part of this change that is expressible through the Excel object
# model is the rename of the first worksheet from "SignalFeed_v2" to
# "Primary", plus that sheet becoming the active tab/selection (as seen
# in the saved sheetView/workbookView state: tabSelected moves from the
# last sheet onto the renamed one, and its selection becomes D15).

$wb = $excel.ActiveWorkbook

# Grab the first sheet defensively by its current (pre-edit) name so this
# still works even if sheet order/index ever changes.
$primary = $null
foreach ($sh in $wb.Worksheets) {
    if ($sh.Name -eq "SignalFeed_v2") {
        $primary = $sh
        break
    }
}
if ($primary -eq $null) {
    $primary = $wb.Worksheets.Item(1)
}

# Rename "SignalFeed_v2" -> "Primary"
$primary.Name = "Primary"

# Activate it (becomes the selected/visible tab, so tabSelected="1" moves
# here and off of whatever sheet previously had it) and move the
# in-sheet selection to D15.
$primary.Activate()
$primary.Range("D15").Select()
